$d = $word.ActiveDocument

# --- Edit 1: append " = $1.109" after "Total = $0.2218 /PCB" ---
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("Total = `$0.2218 /PCB", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $rng1.Collapse(0)
    $rng1.InsertAfter(" = `$1.109")
}

# --- Edit 2: append two sentences about input voltage after the LED sentence ---
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("an LED turns on. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $rng2.Collapse(0)
    $rng2.InsertAfter("It also receives an input voltage from a 5V Lithium-Ion Battery. It Interfaces with every other circuit on the board via its output 3.3V as their input voltage/power supply.")
}
